$d = $word.ActiveDocument

# 1. Replace the empty paragraph right after the title with a new paragraph
#    containing the date text "11/14/2024".
$dateParagraph = $d.Paragraphs.Item(3)
$dateParagraph.Range.InsertBefore("11/14/2024")

# 2. Mark the runs that hold the two inline drawings as "no proofing" so
#    Word stamps <w:rPr><w:noProof/></w:rPr> on them (screenshots pasted
#    in from outside the document).
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = $true
}
